# Apply updates described by the diff:
# - Update several odds values in row 2
# - Delete row 5 entirely (shrinks used range from A1:BD5 to A1:BD4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values in row 2
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 1.83
$ws.Range("K2").Value = 2.6
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.95
$ws.Range("W2").Value = 8.5
$ws.Range("X2").Value = 7.5
$ws.Range("AG2").Value = 201
$ws.Range("AH2").Value = 21
$ws.Range("AN2").Value = 3.5
$ws.Range("AP2").Value = 15
$ws.Range("AT2").Value = 3.5
$ws.Range("BB2").Value = 201
$ws.Range("BC2").Value = 501

# Delete row 5 (the Uruguay / Racing Montevideo - Nacional match)
$ws.Rows("5:5").Delete()
